{"js": "// Replace the date line and each \"NNN\u00d7N=\" expression with its updated\n// value. All original values are unique within the document, so a\n// search-and-replace keyed on the old text is safe and unambiguous.\nconst replacements = [\n  [\"2025-03-01 Saturday\", \"2025-03-02 Sunday\"],\n  [\"981\u00d78=\", \"822\u00d72=\"],\n  [\"353\u00d74=\", \"956\u00d77=\"],\n  [\"733\u00d78=\", \"966\u00d75=\"],\n  [\"241\u00d74=\", \"736\u00d77=\"],\n  [\"280\u00d73=\", \"503\u00d72=\"],\n  [\"626\u00d78=\", \"210\u00d73=\"],\n  [\"408\u00d75=\", \"830\u00d77=\"],\n  [\"224\u00d73=\", \"565\u00d74=\"],\n  [\"674\u00d77=\", \"589\u00d77=\"],\n  [\"906\u00d75=\", \"262\u00d79=\"],\n  [\"658\u00d72=\", \"284\u00d78=\"],\n  [\"193\u00d74=\", \"253\u00d72=\"],\n  [\"733\u00d74=\", \"433\u00d74=\"],\n  [\"512\u00d72=\", \"918\u00d78=\"],\n  [\"298\u00d79=\", \"623\u00d76=\"],\n  [\"802\u00d78=\", \"797\u00d79=\"],\n  [\"678\u00d78=\", \"363\u00d76=\"],\n  [\"602\u00d76=\", \"250\u00d75=\"],\n  [\"457\u00d77=\", \"506\u00d78=\"],\n  [\"221\u00d77=\", \"797\u00d76=\"],\n  [\"328\u00d78=\", \"726\u00d75=\"],\n  [\"102\u00d72=\", \"920\u00d77=\"],\n  [\"642\u00d74=\", \"457\u00d72=\"],\n  [\"705\u00d79=\", \"992\u00d76=\"],\n  [\"815\u00d75=\", \"630\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each \"NNN\u00d7N=\" expression with its updated\n# value. All original values are unique within the document, so a\n# Find/Replace keyed on the old text is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-01 Saturday\", \"2025-03-02 Sunday\"),\n    @(\"981\u00d78=\", \"822\u00d72=\"),\n    @(\"353\u00d74=\", \"956\u00d77=\"),\n    @(\"733\u00d78=\", \"966\u00d75=\"),\n    @(\"241\u00d74=\", \"736\u00d77=\"),\n    @(\"280\u00d73=\", \"503\u00d72=\"),\n    @(\"626\u00d78=\", \"210\u00d73=\"),\n    @(\"408\u00d75=\", \"830\u00d77=\"),\n    @(\"224\u00d73=\", \"565\u00d74=\"),\n    @(\"674\u00d77=\", \"589\u00d77=\"),\n    @(\"906\u00d75=\", \"262\u00d79=\"),\n    @(\"658\u00d72=\", \"284\u00d78=\"),\n    @(\"193\u00d74=\", \"253\u00d72=\"),\n    @(\"733\u00d74=\", \"433\u00d74=\"),\n    @(\"512\u00d72=\", \"918\u00d78=\"),\n    @(\"298\u00d79=\", \"623\u00d76=\"),\n    @(\"802\u00d78=\", \"797\u00d79=\"),\n    @(\"678\u00d78=\", \"363\u00d76=\"),\n    @(\"602\u00d76=\", \"250\u00d75=\"),\n    @(\"457\u00d77=\", \"506\u00d78=\"),\n    @(\"221\u00d77=\", \"797\u00d76=\"),\n    @(\"328\u00d78=\", \"726\u00d75=\"),\n    @(\"102\u00d72=\", \"920\u00d77=\"),\n    @(\"642\u00d74=\", \"457\u00d72=\"),\n    @(\"705\u00d79=\", \"992\u00d76=\"),\n    @(\"815\u00d75=\", \"630\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]2, [ref]$false, [ref]$newText, [ref]2) | Out-Null\n}\n\n$d.Save()\n"}
